$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.614.50'
$ws.Range("E2").Value = '  +0.25%  '
$ws.Range("D3").Value = '3.697.38'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '676.58'
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.84'
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  +0.52%  '
$ws.Range("E9").Value = '  +1.22%  '
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("E11").Value = '  +1.40%  '
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.60'
$ws.Range("E13").Value = '  +0.53%  '
$ws.Range("D14").Value = '3.707.66'
$ws.Range("E14").Value = '  +0.60%  '
$ws.Range("D15").Value = '69.621.22'
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("E16").Value = '  +2.01%  '
$ws.Range("E17").Value = '  +1.17%  '
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '470.98'
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.82'
$ws.Range("E20").Value = '  -1.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.652'
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '80.58'
$ws.Range("E22").Value = '  +1.15%  '
$ws.Range("D23").Value = '3.843.16'
$ws.Range("E23").Value = '  +0.39%  '
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("E25").Value = '  +3.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.90'
$ws.Range("E26").Value = '  -0.32%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.12'
$ws.Range("E27").Value = '  -0.63%  '
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.75'
$ws.Range("E29").Value = '  +1.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.02'
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.60'
$ws.Range("E31").Value = '  +0.64%  '
$ws.Range("E32").Value = '  +0.41%  '
$ws.Range("E33").Value = '  +0.41%  '
$ws.Range("D34").Value = '3.686.20'
$ws.Range("E34").Value = '  +0.77%  '
$ws.Range("E35").Value = '  +0.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.49'
$ws.Range("E36").Value = '  +4.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.23'
$ws.Range("E37").Value = '  +1.48%  '
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0902'
$ws.Range("E41").Value = '  +0.13%  '
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '166.97'
$ws.Range("E43").Value = '  +1.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '47.00'
$ws.Range("E44").Value = '  -0.66%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.77'
$ws.Range("E45").Value = '  +2.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '28.22'
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("E47").Value = '  +1.15%  '
$ws.Range("E48").Value = '  +0.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.10'
$ws.Range("E49").Value = '  -1.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.87'
$ws.Range("E50").Value = '  +0.78%  '
$ws.Range("E51").Value = '  +1.82%  '
